# cryptos.xlsx - "Updated symbol list" refresh (GitHub Actions scraper run).
#
# The sheet stores every value (including the numeric-looking "Price" column D)
# as text, so we must avoid Excel's automatic text->number coercion when
# writing new prices. Set-TextValue forces the cell to Text format for the
# write, then resets the style back to Normal so no residual number format is
# left behind on the cell (matches the source file, which carries no special
# style on these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Price (column D) refresh for existing rows ---------------------------
Set-TextValue "D2"  "245.22"
Set-TextValue "D3"  "23.93"
Set-TextValue "D4"  "5.199"
Set-TextValue "D5"  "0.05742"
Set-TextValue "D6"  "6.459"
Set-TextValue "D7"  "3.210"
Set-TextValue "D8"  "0.8140"
Set-TextValue "D9"  "0.8683"
Set-TextValue "D10" "0.1378"
Set-TextValue "D11" "0.06948"
Set-TextValue "D12" "0.03163"
Set-TextValue "D13" "0.03016"
Set-TextValue "D14" "0.09332"
Set-TextValue "D15" "3.816"
Set-TextValue "D16" "0.001524"
Set-TextValue "D17" "0.04714"
Set-TextValue "D18" "0.0005976"
Set-TextValue "D19" "0.006260"
Set-TextValue "D20" "0.001234"
Set-TextValue "D21" "0.004110"
Set-TextValue "D22" "0.00008695"
Set-TextValue "D23" "3.585"
Set-TextValue "D24" "2.145"
Set-TextValue "D26" "0.1329"
Set-TextValue "D27" "0.0002327"
Set-TextValue "D40" "0.03718"

# --- Rows 41-43 reshuffle: KickToken / BKEXToken / CEJI rotate position ----
# Row 41: was KickToken -> now BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1052"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42: was BKEXToken -> now CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.002309"
$ws.Range("E42").Value = "41CEJICEJIWorstin24h"

# Row 43: was CEJI -> now KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.006226"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- Remaining trailing updates --------------------------------------------
Set-TextValue "D44" "0.007459"
Set-TextValue "D45" "0.00005374"
Set-TextValue "D47" "0.4397"
Set-TextValue "D48" "0.002244"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
